$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric D/E/F values (relative abundance recomputed per-column) ---

# Row 2
$ws.Range("E2").Value = 0.08015624140397205

# Row 3
$ws.Range("F3").Value = 0.4732762888056007

# Row 4 (new zeros added)
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

# Row 5 (new zeros added)
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# Row 6
$ws.Range("F6").Value = 0.02574075383636235

# Row 7
$ws.Range("D7").Value = 0.02901290799401867

# Row 8
$ws.Range("E8").Value = 0.001283673506812639

# Row 9 (new zeros added)
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0

# Row 10
$ws.Range("E10").Value = 0.003062478223395866

# Row 11
$ws.Range("D11").Value = 0.04035681751774634

# Row 12 (new zeros added)
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0

# Row 13
$ws.Range("D13").Value = 0.005465701861432426

# Row 14
$ws.Range("E14").Value = 0.001760466523628762
$ws.Range("F14").Value = 0.0411710628668411

# Row 15
$ws.Range("E15").Value = 0.001925510260218958

# Row 16
$ws.Range("E16").Value = 0.0120665309640388

# Row 17
$ws.Range("E17").Value = 0.0458821587720746

# Row 18 (new zeros added)
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0

# Row 19
$ws.Range("D19").Value = 0.09575290902528316
$ws.Range("E19").Value = 0.003704314976802186

# Row 20 (new zeros added)
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0

# Row 21
$ws.Range("E21").Value = 0.0008618950688599145

# Row 22
$ws.Range("D22").Value = 0.01014076760454444
$ws.Range("E22").Value = 0.1729474977535713
$ws.Range("F22").Value = 0.2076939396082314

# Row 23 (new zeros added)
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0

# Row 24 (new zeros added)
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0

# Row 25
$ws.Range("D25").Value = 0.6262353689348756
$ws.Range("E25").Value = 0.01452384881993728

# Row 26
$ws.Range("E26").Value = 0.0003117492802259266

# Row 27
$ws.Range("D27").Value = 0.001598459978343446
$ws.Range("E27").Value = 0.02416973831398654

# Row 28
$ws.Range("D28").Value = 0.03330984341967309
$ws.Range("E28").Value = 0.01454218701289175
$ws.Range("F28").Value = 0.0004808712255144615

# Row 29
$ws.Range("D29").Value = 0.002165655454529829
$ws.Range("E29").Value = 0.03731822266233886
$ws.Range("F29").Value = 0.0000424298140159819

# Row 30
$ws.Range("E30").Value = 0.006583411270653389

# Row 31
$ws.Range("D31").Value = 0.1302143311390316
$ws.Range("E31").Value = 0.02070381984559241

# Row 32
$ws.Range("E32").Value = 0.003172507381122664

# Row 33
$ws.Range("D33").Value = 0.001289080627696327

# Row 34
$ws.Range("D34").Value = 0.00006875096681047077
$ws.Range("E34").Value = 0.01709119583356256

# Row 35
$ws.Range("D35").Value = 0.002715663189013596
$ws.Range("E35").Value = 0.491188498285379
$ws.Range("F35").Value = 0.2480588360087688

# Row 36
$ws.Range("E36").Value = 0.02099723093286388

# Row 37
$ws.Range("E37").Value = 0.001485393629311768

# Row 38 (new zeros added)
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0

# --- Rows 39 & 40: content swapped (species rows reordered) and values recomputed ---
# Row 39 becomes the "Unassigned" entry
$ws.Range("A39").Value = "Unassigned"
$ws.Range("B39").Value = "Unassigned"
$ws.Range("C39").Value = "Unassigned"
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0.005464781500430948
$ws.Range("F39").Value = 0.003380241849939891

# Row 40 becomes the "Urophycis sp" entry
$ws.Range("A40").Value = "Urophycis sp"
$ws.Range("B40").Value = "Red White or Spotted hake"
$ws.Range("C40").Value = "Teleost Fish"
$ws.Range("D40").Value = 0.02167374228700091
$ws.Range("E40").Value = 0.01879664777832792
$ws.Range("F40").Value = 0.000155575984725267

$wb.Save()
